$d = $word.ActiveDocument

# 1) Vertically center the content of the "Payment" cell (table 4, row 1, col 1).
$paymentTable = $d.Tables.Item(4)
$paymentCell = $paymentTable.Cell(1, 1)
$paymentCell.VerticalAlignment = 1

# 2) Remove the standalone "Pembayaran akan dilakukan dengan" paragraph; its wording is
#    folded into the paragraph that follows it (which used to read
#    "mengurangi CazhBOX lembaga Anda.").
$paraCount = $d.Content.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $d.Content.Paragraphs.Item($i)
    if ($para.Range.Text -eq "Pembayaran akan dilakukan dengan" + [char]13) {
        $para.Range.Delete()
        break
    }
}

# 3) Rewrite the remaining paragraph so it reads:
#    "Pembayaran akan dilakukan dengan ${metode}."
$rng = $d.Content
$rng.Find.Execute("mengurangi ", $false, $false, $false, $false, $false, $true, 1, $false, "Pembayaran akan dilakukan dengan ", 2)

$cazh = $d.Content
$cazh.Find.Execute("CazhBOX lembaga Anda.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($cazh.Find.Found) {
    $cazh.Text = '${metode}.'
}
